$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'26.419.67"
$ws.Range('E2').Value = '  +1.29%  '
$ws.Range('D3').Value = "'1.676.84"
$ws.Range('E3').Value = '  +2.41%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = "'217.03"
$ws.Range('E5').Value = '  +1.43%  '
$ws.Range('D6').Value = "'0.5301"
$ws.Range('E6').Value = '  +0.88%  '
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('E8').Value = '  +3.85%  '
$ws.Range('D9').Value = "'0.06402"
$ws.Range('E9').Value = '  +1.46%  '
$ws.Range('D10').Value = "'21.73"
$ws.Range('E10').Value = '  +5.01%  '
$ws.Range('D11').Value = "'0.07819"
$ws.Range('E11').Value = '  +2.18%  '
$ws.Range('D12').Value = "'1.688.84"
$ws.Range('E12').Value = '  +3.81%  '
$ws.Range('D13').Value = "'4.511"
$ws.Range('E13').Value = '  +1.99%  '
$ws.Range('D14').Value = "'0.5567"
$ws.Range('E14').Value = '  +1.26%  '
$ws.Range('D15').Value = "'0.0₅8342"
$ws.Range('E15').Value = '  +2.15%  '
$ws.Range('D16').Value = "'65.63"
$ws.Range('E16').Value = '  +0.88%  '
$ws.Range('D17').Value = "'26.470.53"
$ws.Range('E17').Value = '  +1.52%  '
$ws.Range('E18').Value = '  -0.08%  '
$ws.Range('D19').Value = "'4.731"
$ws.Range('E19').Value = '  +0.84%  '
$ws.Range('D20').Value = "'193.54"
$ws.Range('E20').Value = '  +2.94%  '
$ws.Range('E21').Value = '  +1.45%  '
$ws.Range('D22').Value = "'6.343"
$ws.Range('E22').Value = '  +2.92%  '
$ws.Range('E23').Value = '  -0.03%  '
$ws.Range('D24').Value = "'142.41"
$ws.Range('E24').Value = '  -2.43%  '
$ws.Range('D25').Value = "'0.1288"
$ws.Range('E25').Value = '  +5.83%  '
$ws.Range('D26').Value = "'7.404"
$ws.Range('E26').Value = '  -0.13%  '
$ws.Range('D27').Value = "'16.25"
$ws.Range('E27').Value = '  +2.58%  '
$ws.Range('E28').Value = '  +1.86%  '
$ws.Range('D29').Value = "'0.06345"
$ws.Range('E29').Value = '  +5.41%  '
$ws.Range('D30').Value = "'1.274"
$ws.Range('E30').Value = '  +1.37%  '
$ws.Range('D31').Value = "'3.622"
$ws.Range('E31').Value = '  +5.18%  '
$ws.Range('D32').Value = "'3.451"
$ws.Range('E32').Value = '  +1.17%  '
$ws.Range('E33').Value = '  +2.25%  '
$ws.Range('E34').Value = '  +2.19%  '
$ws.Range('D35').Value = "'0.6207"
$ws.Range('E35').Value = '  +8.15%  '
$ws.Range('D36').Value = "'2.428"
$ws.Range('E36').Value = '  +1.31%  '
$ws.Range('D37').Value = "'2.781"
$ws.Range('E37').Value = '  +0.53%  '
$ws.Range('D38').Value = "'6.162"
$ws.Range('E38').Value = '  +7.55%  '
$ws.Range('E39').Value = '  +0.85%  '
$ws.Range('D40').Value = "'1.084.10"
$ws.Range('E40').Value = '  +4.37%  '
$ws.Range('D41').Value = "'0.8652"
$ws.Range('E41').Value = '  +1.00%  '
$ws.Range('D42').Value = "'1.0000"
$ws.Range('E42').Value = '  -0.14%  '
$ws.Range('D43').Value = "'100.38"
$ws.Range('E43').Value = '  -0.30%  '
$ws.Range('D44').Value = "'1.822.60"
$ws.Range('E44').Value = '  +1.96%  '
$ws.Range('D45').Value = "'57.32"
$ws.Range('E45').Value = '  +3.06%  '
$ws.Range('D46').Value = "'8.155"
$ws.Range('E46').Value = '  +1.32%  '
$ws.Range('D47').Value = "'1.002"
$ws.Range('E47').Value = '  -0.07%  '
$ws.Range('E48').Value = '  -4.78%  '
$ws.Range('D49').Value = "'0.05209"
$ws.Range('B50').Value = 'Aptos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D50').Value = "'6.034"
$ws.Range('E50').Value = '  +2.01%  '
$ws.Range('B51').Value = 'RenderToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D51').Value = "'1.467"
$ws.Range('E51').Value = '  +5.39%  '
